$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix two typos in existing cells ---
$ws.Range("H2").Value = "Sam want learn something"
$ws.Range("H3").Value = "Mart want spent time with friends with passive developing"

# --- New "Vlad / disabled person" user entry on row 4 ---
$ws.Range("G4").Value = "Vlad"
$ws.Range("H4").Value = "disabled person"

# --- Remove the old "ASM want view appropriate context..." text from C22 ---
$ws.Range("C22").Clear()

# --- New BackLog section content (rows 28-46) ---
$ws.Range("C28").Value = "ASM want to use app with mobile"
$ws.Range("E28").Value = "1. Adaptive design"

$ws.Range("C29").Value = "ASM want differntiate each other in multiptiplayer game"
$ws.Range("E29").Value = "1. Avatars"

$ws.Range("E30").Value = "2. Names (Nicknames)"

$ws.Range("E31").Value = "3. Identificators"

$ws.Range("E32").Value = "4. User colors for hero and for enemy (question visualisation)"
$ws.Rows.Item(32).RowHeight = 23.85

$ws.Range("C33").Value = "ASM want to have different variant of game formats"
$ws.Range("E33").Value = "1.  Crossword"

$ws.Range("E34").Value = "2. True-false"

$ws.Range("E35").Value = "3. Quiz (choose answer from list)"

$ws.Range("E36").Value = "4. Quiz (write answer manually)"

$ws.Range("C37").Value = "ASM want to have possibility to communicate with each other"
$ws.Range("E37").Value = "1. Add description to account with contacts (mail, telegram etc.)"
$ws.Rows.Item(37).RowHeight = 23.85

$ws.Range("E38").Value = "2. Voice chat"

$ws.Range("E39").Value = "3. Text chat"

$ws.Range("E40").Value = "4. Add Loud Phrase (user pres on button and all team members see phrase)"
$ws.Rows.Item(40).RowHeight = 35.05

$ws.Range("C41").Value = "S want to analyse progress and use it for learning"
$ws.Range("E41").Value = "1. Statistics"

$ws.Range("E42").Value = "2. Repeat mode"

$ws.Range("E43").Value = "3. Compare with other players"

$ws.Range("C44").Value = "VS want to have accessible interface"
$ws.Range("E44").Value = "1. Voice interface"

$ws.Range("E45").Value = "2. Big font"

$ws.Range("E46").Value = "3. Keyboard first"

# --- Update selection to match the area just edited ---
$ws.Range("C28:E46").Select()
